# Apply weekly update: insert a new record row at row 63 (pushing existing
# rows 63-72 down to 64-73) and populate the new row with the latest
# Fruta/Mercado Mayorista Lo Valledor de Santiago - Coco price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63; this shifts rows 63:72 down to 64:73
# and carries formatting from the row above (keeping the D-column date style).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Range("A63").Value = 6
$ws.Range("B63").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44776
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100108
$ws.Range("H63").Value = "Tropicales y subtropicales"
$ws.Range("I63").Value = 100108007
$ws.Range("J63").Value = "Coco"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 150
$ws.Range("N63").Value = 23000
$ws.Range("O63").Value = 24000
$ws.Range("P63").Value = 23500
$ws.Range("Q63").Value = "`$/malla 20 unidades"
$ws.Range("R63").Value = "Perú"
$ws.Range("S63").Value = 1175
$ws.Range("T63").Value = 20
